$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'332.07"
$ws.Range("E2").Value = "'2.26%"
$ws.Range("D3").Value = "'41.10"
$ws.Range("E3").Value = "'3.18%"
$ws.Range("D4").Value = "'5.732"
$ws.Range("E4").Value = "'-1.85%"
$ws.Range("D5").Value = "'0.08200"
$ws.Range("E5").Value = "'2.73%"
$ws.Range("D6").Value = "'2.044"
$ws.Range("E6").Value = "'6.65%"
$ws.Range("D7").Value = "'8.747"
$ws.Range("E7").Value = "'0.60%"
$ws.Range("D8").Value = "'4.513"
$ws.Range("E8").Value = "'-1.33%"
$ws.Range("E9").Value = "'1.75%"
$ws.Range("D10").Value = "'0.9217"
$ws.Range("E10").Value = "'-1.83%"
$ws.Range("D11").Value = "'0.1242"
$ws.Range("E11").Value = "'-1.46%"
$ws.Range("D12").Value = "'0.1952"
$ws.Range("E12").Value = "'-0.34%"
$ws.Range("D13").Value = "'8.323"
$ws.Range("E13").Value = "'-5.34%"
$ws.Range("D14").Value = "'0.09424"
$ws.Range("E14").Value = "'2.70%"
$ws.Range("D15").Value = "'0.03634"
$ws.Range("E15").Value = "'1.87%"
$ws.Range("E16").Value = "'9.64%"
$ws.Range("D17").Value = "'0.001301"
$ws.Range("E17").Value = "'-0.59%"
$ws.Range("D18").Value = "'0.006214"
$ws.Range("E18").Value = "'1.02%"
$ws.Range("D19").Value = "'3.388"
$ws.Range("E19").Value = "'1.17%"
$ws.Range("E20").Value = "'-1.16%"
$ws.Range("D21").Value = "'0.1417"
$ws.Range("E21").Value = "'-1.11%"
$ws.Range("D22").Value = "'0.2650"
$ws.Range("E22").Value = "'9.68%"
$ws.Range("D23").Value = "'0.04426"
$ws.Range("E23").Value = "'-0.58%"
$ws.Range("E24").Value = "'-0.04%"
$ws.Range("D25").Value = "'0.004316"
$ws.Range("E25").Value = "'-1.86%"
$ws.Range("E26").Value = "'8.44%"
$ws.Range("D39").Value = "'0.02775"
$ws.Range("E39").Value = "'14.73%"
$ws.Range("D40").Value = "'0.05516"
$ws.Range("E40").Value = "'5.29%"
$ws.Range("D41").Value = "'0.007617"
$ws.Range("E41").Value = "'2.21%"
$ws.Range("D42").Value = "'0.009952"
$ws.Range("E42").Value = "'14.41%"
$ws.Range("D43").Value = "'0.1423"
$ws.Range("E43").Value = "'0.87%"
$ws.Range("D44").Value = "'0.002120"
$ws.Range("E44").Value = "'-0.37%"
$ws.Range("D45").Value = "'0.01195"
$ws.Range("E45").Value = "'13.57%"
$ws.Range("D46").Value = "'0.00006754"
$ws.Range("E46").Value = "'-1.05%"
$ws.Range("E47").Value = "'-0.38%"
$ws.Range("D48").Value = "'0.002279"
$ws.Range("E48").Value = "'59.86%"
$ws.Range("D49").Value = "'0.002993"
$ws.Range("E49").Value = "'4.00%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.38%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.38%"
